# Insert a new weekly price record as row 321 in the "Perejil" sheet.
# This shifts all existing rows from 321..373 down to 322..374 and
# updates the worksheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 321 (pushes 321:373 -> 322:374)
$ws.Rows.Item(321).Insert()

# Populate the new row 321 with the new record's data.
$ws.Cells.Item(321, 1).Value = 9
$ws.Cells.Item(321, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(321, 3).Value = "Metropolitana"
$ws.Cells.Item(321, 4).Value = 44694
$ws.Cells.Item(321, 5).Value = 13
$ws.Cells.Item(321, 6).Value = 100112044
$ws.Cells.Item(321, 7).Value = "Perejil"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 79
$ws.Cells.Item(321, 11).Value = 9000
$ws.Cells.Item(321, 12).Value = 10000
$ws.Cells.Item(321, 13).Value = 9494
$ws.Cells.Item(321, 14).Value = "$/docena de atados"
$ws.Cells.Item(321, 15).Value = "Región Metropolitana"
$ws.Cells.Item(321, 16).Value = 3165
$ws.Cells.Item(321, 17).Value = 3
$ws.Cells.Item(321, 18).Value = "Hortaliza"

# Keep the date column formatted the same as the surrounding rows (row
# insert already copies this from the row above, but set explicitly to
# be safe).
$ws.Cells.Item(321, 4).NumberFormat = $ws.Cells.Item(322, 4).NumberFormat
